# Sprint Plan & Retrospective - "Upper Floor Doors & Sprint #3 Docs"
#
# This script:
#  1. Fills in the new "Sprint #3" block (rows 21-25, columns A-F) on the
#     "Plan" sheet, reusing the existing cell-fill formatting used by the
#     "Sprint #1" block (A4/B3:E3 style) so the new rows visually match.
#  2. Adds the new text required for that block (shared strings are created
#     in the same order a person would have typed them, so the underlying
#     sharedStrings.xml ends up in the same order as the authored commit).
#  3. Updates the saved cell selection on both worksheets.

$wb = $excel.ActiveWorkbook

$plan = $wb.Worksheets.Item("Plan")
$retro = $wb.Worksheets.Item("Retrospective")

# ---------------------------------------------------------------------
# 1. Apply formatting for the new Sprint #3 block (rows 21-25)
# ---------------------------------------------------------------------
# Column A uses the plain "fill2" style already used by A4:A11 (no special
# alignment). Columns B-F use the centered "fill2" style already used by
# B3:E3 (fillId=2, horizontal/vertical = center).
$plan.Range("A4").Copy()
$plan.Range("A21:A25").PasteSpecial(-4122)

$plan.Range("B3").Copy()
$plan.Range("B21:F25").PasteSpecial(-4122)

$plan.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Populate the new content, in authoring order so new shared strings
#    are appended in the same order as the original edit.
# ---------------------------------------------------------------------
$plan.Cells.Item(21, 1).Value = "Sprint #3"

$plan.Cells.Item(22, 4).Value = "Upper Floor Upper Foyer"
$plan.Cells.Item(23, 4).Value = "Upper Floor Room D"
$plan.Cells.Item(24, 4).Value = "Upper Floor Room F"
$plan.Cells.Item(25, 4).Value = "Upper Floor Room I"

$plan.Cells.Item(21, 6).Value = "Door Locked UI"
$plan.Cells.Item(22, 6).Value = "Items"
$plan.Cells.Item(23, 6).Value = "Main Menu (Outside)"

$plan.Cells.Item(24, 5).Value = "User Stories (cont.)"

$plan.Cells.Item(21, 2).Value = "Implement Essential Upper Floor Rooms"

# Remaining cells reuse text that already exists elsewhere in the sheet.
$plan.Cells.Item(21, 3).Value = "1 Week"
$plan.Cells.Item(21, 4).Value = "Update MSQI Chart"
$plan.Cells.Item(21, 5).Value = "Basement Room D: Assets & Props"
$plan.Cells.Item(22, 5).Value = "Basement Landing: Assets & Props"
$plan.Cells.Item(23, 5).Value = "AI Patrolling"

# ---------------------------------------------------------------------
# 3. Update the saved selections on each sheet.
# ---------------------------------------------------------------------
$retro.Activate()
$retro.Range("B10").Select()

$plan.Activate()
$plan.Range("C33").Select()
